$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2-46, 49-51: update Price (D) and Volume(1h) (E) values ---
# Price cells are set via a leading apostrophe to force literal text
# (matching the original inline-string cell type) and the cell style
# is reset to Normal afterwards so no stray quote-prefix number format
# sticks around on cells that did not have one before.

$ws.Range("D2").Value = "'67.965.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.64%  "

$ws.Range("D3").Value = "'3.801.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.35%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'594.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.06%  "

$ws.Range("D6").Value = "'171.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.93%  "

$ws.Range("D7").Value = "'3.794.97"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.07%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "'0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("E10").Value = "  -3.47%  "

$ws.Range("D11").Value = "'6.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("D12").Value = "'0.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.28%  "

$ws.Range("D13").Value = "'38.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.24%  "

$ws.Range("E14").Value = "  -3.94%  "

$ws.Range("D15").Value = "'4.434.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.24%  "

$ws.Range("D16").Value = "'3.794.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.46%  "

$ws.Range("D17").Value = "'68.251.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.23%  "

$ws.Range("E18").Value = "  -4.38%  "

$ws.Range("E19").Value = "  -4.07%  "

$ws.Range("D20").Value = "'16.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("D21").Value = "'488.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.26%  "

$ws.Range("D22").Value = "'9.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.47%  "

$ws.Range("D23").Value = "'0.739"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.09%  "

$ws.Range("D24").Value = "'86.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("E25").Value = "  -5.49%  "

$ws.Range("D26").Value = "'0.0000138"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.72%  "

$ws.Range("D27").Value = "'12.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.54%  "

$ws.Range("E28").Value = "  -7.80%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").Value = "'32.53"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = "'7.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.94%  "

$ws.Range("E34").Value = "  -3.01%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("D36").Value = "'1.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.38%  "

$ws.Range("D37").Value = "'5.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.23%  "

$ws.Range("E38").Value = "  -2.31%  "

$ws.Range("D39").Value = "'0.326"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.22%  "

$ws.Range("D40").Value = "'450.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.98%  "

$ws.Range("D41").Value = "'49.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("E42").Value = "  -2.34%  "

$ws.Range("D43").Value = "'2.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.89%  "

$ws.Range("D44").Value = "'8.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("D45").Value = "'41.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.69%  "

$ws.Range("D46").Value = "'2.862.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.88%  "

# --- Rows 47-48: VeChain and USDe swap places (row 47 becomes USDe,
#     row 48 becomes VeChain), with freshly-updated Volume(1h) values ---
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0353"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("D49").Value = "'137.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").Value = "'26.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.97%  "

$ws.Range("D51").Value = "'23.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.78%  "

